# Auto commit at 2026-02-13  8:01:10.63
#
# Update the Metrics sheet's raw input figures for this period and clear
# the one-off manual adjustment formulas on the "today" sheet (B3:B6).
# Every other touched cell (today!B11:F22, today!A1's TODAY()-1 cache,
# etc.) is a live formula and will simply recompute from these inputs.

$wb = $excel.ActiveWorkbook

# --- Metrics sheet: refresh B2:B13 with the latest figures -----------------
$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value  = 190429.38
$metrics.Range("B3").Value  = 172998.75
$metrics.Range("B4").Value  = 65512.72
$metrics.Range("B5").Value  = 7626
$metrics.Range("B6").Value  = 771065.15
$metrics.Range("B7").Value  = 625793.19999999995
$metrics.Range("B8").Value  = 229947.28
$metrics.Range("B9").Value  = 31102
$metrics.Range("B10").Value = 34872316.870000005
$metrics.Range("B11").Value = 32671785.989999998
$metrics.Range("B12").Value = 12175761.140000001
$metrics.Range("B13").Value = 1349009

$metrics.Activate()
$metrics.Range("D13").Select()

# --- today sheet: clear the stale manual-adjustment formulas ---------------
$today = $wb.Worksheets.Item("today")

$today.Range("B3:B6").ClearContents()

$today.Activate()
$today.Range("H21").Select()
